$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 188.4
$ws.Range("I2").Value = 188.4
$ws.Range("K2").Value = 188.4
$ws.Range("M2").Value = -75.40000000000001
$ws.Range("H137").Value = 842.4286
$ws.Range("I137").Value = 779.6
$ws.Range("K137").Value = 2338.8
$ws.Range("M137").Value = 211.1999999999998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 284879.5
$ws.Range("J69").Value = 284879.5
$ws.Range("L69").Value = 284879.5
$ws.Range("N69").Value = -286377.5
$ws.Range("H72").Value = 284879.5
$ws.Range("J72").Value = 284879.5
$ws.Range("L72").Value = 854638.5
$ws.Range("N72").Value = -862126.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 781.625
$ws.Range("I5").Value = 814.6667
$ws.Range("J5").Value = 761.8
$ws.Range("K5").Value = 814.6667
$ws.Range("L5").Value = 761.8
$ws.Range("M5").Value = -701.6667
$ws.Range("N5").Value = -987.8
$ws.Range("H70").Value = 347500
$ws.Range("J70").Value = 347500
$ws.Range("L70").Value = 347500
$ws.Range("N70").Value = -348086
$ws.Range("H73").Value = 347500
$ws.Range("J73").Value = 347500
$ws.Range("L73").Value = 347500
$ws.Range("N73").Value = -349528
$ws.Range("H86").Value = 13309.667
$ws.Range("I86").Value = 18133.285
$ws.Range("J86").Value = 6556.6
$ws.Range("K86").Value = 18133.285
$ws.Range("L86").Value = 6556.6
$ws.Range("M86").Value = -17010.285
$ws.Range("N86").Value = -8802.6
$ws.Range("H89").Value = 13309.667
$ws.Range("I89").Value = 18133.285
$ws.Range("J89").Value = 6556.6
$ws.Range("K89").Value = 90666.425
$ws.Range("L89").Value = 32783
$ws.Range("M89").Value = -85050.425
$ws.Range("N89").Value = -44015
$ws.Range("H105").Value = 2184.6
$ws.Range("I105").Value = 1530
$ws.Range("K105").Value = 1530
$ws.Range("M105").Value = 217
$ws.Range("H134").Value = 1995
$ws.Range("I134").Value = 1995
$ws.Range("K134").Value = 5985
$ws.Range("M134").Value = -3450
$ws.Range("H135").Value = 61873
$ws.Range("J135").Value = 61873
$ws.Range("L135").Value = 61873
$ws.Range("N135").Value = -72013

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 3400
$ws.Range("I36").Value = 3400
$ws.Range("K36").Value = 3400
$ws.Range("M36").Value = -3012
$ws.Range("H40").Value = 3400
$ws.Range("I40").Value = 3400
$ws.Range("K40").Value = 3400
$ws.Range("M40").Value = -3240
$ws.Range("H107").Value = 795.7143
$ws.Range("I107").Value = 575
$ws.Range("J107").Value = 1347.5
$ws.Range("K107").Value = 575
$ws.Range("L107").Value = 1347.5
$ws.Range("M107").Value = 1345
$ws.Range("N107").Value = -5187.5
$ws.Range("H122").Value = 1425.2222
$ws.Range("I122").Value = 1447.125
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 4341.375
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -1891.375
$ws.Range("N122").Value = -8650

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2137192.5
$ws.Range("I4").Value = 795740.75
$ws.Range("K4").Value = 2387222.25
$ws.Range("M4").Value = -2387110.25
$ws.Range("H7").Value = 93082.38
$ws.Range("I7").Value = 130030.78
$ws.Range("J7").Value = 9948.5
$ws.Range("K7").Value = 390092.34
$ws.Range("L7").Value = 29845.5
$ws.Range("M7").Value = -389980.34
$ws.Range("N7").Value = -30069.5
$ws.Range("H23").Value = 291.4
$ws.Range("I23").Value = 301
$ws.Range("J23").Value = 281.8
$ws.Range("K23").Value = 903
$ws.Range("L23").Value = 845.4000000000001
$ws.Range("M23").Value = -668
$ws.Range("N23").Value = -1315.4
$ws.Range("H24").Value = 124.5
$ws.Range("I24").Value = 50
$ws.Range("J24").Value = 199
$ws.Range("K24").Value = 150
$ws.Range("L24").Value = 597
$ws.Range("M24").Value = 80
$ws.Range("N24").Value = -1057
$ws.Range("H26").Value = 531.6667
$ws.Range("I26").Value = 408.75
$ws.Range("J26").Value = 777.5
$ws.Range("K26").Value = 1226.25
$ws.Range("L26").Value = 2332.5
$ws.Range("M26").Value = -938.25
$ws.Range("N26").Value = -2908.5
$ws.Range("H34").Value = 3652.4783
$ws.Range("I34").Value = 384
$ws.Range("J34").Value = 4142.75
$ws.Range("K34").Value = 1152
$ws.Range("L34").Value = 12428.25
$ws.Range("M34").Value = -1068
$ws.Range("N34").Value = -12596.25
$ws.Range("H49").Value = 2733
$ws.Range("I49").Value = 2999.5
$ws.Range("J49").Value = 2200
$ws.Range("K49").Value = 8998.5
$ws.Range("L49").Value = 6600
$ws.Range("M49").Value = -8842.5
$ws.Range("N49").Value = -6912

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 8025106
$ws.Range("I14").Value = 11464368
$ws.Range("J14").Value = 161.66667
$ws.Range("K14").Value = 11464368
$ws.Range("L14").Value = 161.66667
$ws.Range("M14").Value = -11464200
$ws.Range("N14").Value = -497.66667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2174.75
$ws.Range("I7").Value = 2166.6667
$ws.Range("J7").Value = 2199
$ws.Range("K7").Value = 2166.6667
$ws.Range("L7").Value = 2199
$ws.Range("M7").Value = -2054.6667
$ws.Range("N7").Value = -2423
$ws.Range("H16").Value = 323
$ws.Range("I16").Value = 323
$ws.Range("K16").Value = 323
$ws.Range("M16").Value = -153
$ws.Range("H40").Value = 4029.5557
$ws.Range("I40").Value = 3897.2856
$ws.Range("J40").Value = 4492.5
$ws.Range("K40").Value = 3897.2856
$ws.Range("L40").Value = 4492.5
$ws.Range("M40").Value = -3761.2856
$ws.Range("N40").Value = -4764.5
$ws.Range("H68").Value = 800
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 800
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -51
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 800
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 4000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -256
$ws.Range("N71").ClearContents()
$ws.Range("H126").Value = 2174.75
$ws.Range("I126").Value = 2166.6667
$ws.Range("J126").Value = 2199
$ws.Range("K126").Value = 6500.000100000001
$ws.Range("L126").Value = 6597
$ws.Range("M126").Value = -4030.000100000001
$ws.Range("N126").Value = -11537
$ws.Range("H132").Value = 4174
$ws.Range("I132").Value = 4199
$ws.Range("K132").Value = 12597
$ws.Range("M132").Value = -10067
$ws.Range("H136").Value = 16615.555
$ws.Range("I136").Value = 13094.667
$ws.Range("K136").Value = 39284.001
$ws.Range("M136").Value = -36734.001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7499
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 9998
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 9998
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -11246
$ws.Range("H65").Value = 7499
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 9998
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 49990
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -56230
$ws.Range("H69").Value = 100271
$ws.Range("J69").Value = 100271
$ws.Range("L69").Value = 100271
$ws.Range("N69").Value = -101769
$ws.Range("H72").Value = 100271
$ws.Range("J72").Value = 100271
$ws.Range("L72").Value = 300813
$ws.Range("N72").Value = -308301
$ws.Range("H136").Value = 57967.668
$ws.Range("I136").Value = 100004
$ws.Range("J136").Value = 36949.5
$ws.Range("K136").Value = 300012
$ws.Range("L136").Value = 110848.5
$ws.Range("M136").Value = -297462
$ws.Range("N136").Value = -115948.5
